$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "Normalized_CanonEx6_7"
$ws.Range("M1").Value = "Normalized_CanonEx7_8"
$ws.Range("N1").Value = "Normalized_ES7"
